# Update paises.xlsx (sheet "Pais") with the latest COVID-19 country stats.
# The source data was re-sorted by total cases, which moves a few countries
# to new rows (Camerun/Azerbaiyan swap around row 72-73; Islas Caimanes
# jumps ahead of Barbados/Liechtenstein around row 167-169) and refreshes
# the "last updated" timestamp plus a handful of per-country counters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Datos actualizados" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 12 de Mayo de 2020 a las 00:05"

# --- Estados Unidos (row 4): refreshed counters, country stays in place ---
$ws.Range("B4").Value = 1381665
$ws.Range("C4").Value = 14027
$ws.Range("E4").Value = 1039925
$ws.Range("G4").Value = 765
$ws.Range("H4").Value = 81552

# --- Argentina (row 56): refreshed counters, country stays in place ---
$ws.Range("B56").Value = 6265
$ws.Range("C56").Value = 231
$ws.Range("E56").Value = 4114
$ws.Range("G56").Value = 9
$ws.Range("H56").Value = 314

# --- Camerun overtakes Azerbaiyan: row 72 becomes Camerun (new data),
#     row 73 becomes Azerbaiyan (its old, unchanged data) ---
$ws.Range("A72").Value = "Camerun"
$ws.Range("B72").Value = 2689
$ws.Range("C72").Value = 110
$ws.Range("D72").Value = 1524
$ws.Range("E72").Value = 1040
$ws.Range("F72").Value = 28
$ws.Range("G72").Value = 11
$ws.Range("H72").Value = 125

$ws.Range("A73").Value = "Azerbaiyan"
$ws.Range("B73").Value = 2589
$ws.Range("C73").Value = 70
$ws.Range("D73").Value = 1680
$ws.Range("E73").Value = 877
$ws.Range("F73").Value = 33
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 32

# --- Islas Caimanes jumps ahead of Barbados and Liechtenstein:
#     row 167 becomes Islas Caimanes (new data), row 168 becomes Barbados
#     (its old data), row 169 becomes Liechtenstein (its old data) ---
$ws.Range("A167").Value = "Islas Caimanes"
$ws.Range("B167").Value = 84
$ws.Range("C167").Value = 3
$ws.Range("D167").Value = 47
$ws.Range("E167").Value = 36
$ws.Range("F167").Value = 3
$ws.Range("G167").Value = 0
$ws.Range("H167").Value = 1

$ws.Range("A168").Value = "Barbados"
$ws.Range("B168").Value = 84
$ws.Range("C168").Value = 0
$ws.Range("D168").Value = 57
$ws.Range("E168").Value = 20
$ws.Range("F168").Value = 4
$ws.Range("G168").Value = 0
$ws.Range("H168").Value = 7

$ws.Range("A169").Value = "Liechtenstein"
$ws.Range("B169").Value = 82
$ws.Range("C169").Value = 0
$ws.Range("D169").Value = 55
$ws.Range("E169").Value = 26
$ws.Range("F169").Value = 0
$ws.Range("G169").Value = 0
$ws.Range("H169").Value = 1
